$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 49385
$ws.Range("E2").Value = 6827
$ws.Range("F2").Value = 6827
$ws.Range("G2").Value = 6830
$ws.Range("H2").Value = 5049
$ws.Range("I2").Value = 5013
$ws.Range("J2").Value = 36
$ws.Range("K2").Value = 56424
$ws.Range("L2").Value = 8593
$ws.Range("M2").Value = 47831
$ws.Range("N2").Value = 46339
$ws.Range("O2").Value = 1492
$ws.Range("P2").Value = 944
$ws.Range("Q2").Value = 4039
$ws.Range("R2").Value = -5636
$ws.Range("S2").Value = -865
$ws.Range("T2").Value = 3916
$ws.Range("U2").Value = 123
$ws.Range("V2").Value = 1720
$ws.Range("W2").Value = 13.82
$ws.Range("X2").Value = 10.22
$ws.Range("Y2").Value = 11.37
$ws.Range("Z2").Value = 9.380000000000001
$ws.Range("AA2").Value = 17.96
$ws.Range("AB2").Value = 4802.51
$ws.Range("AC2").Value = 26565
$ws.Range("AD2").Value = 15.19
$ws.Range("AE2").Value = 262184
$ws.Range("AF2").Value = 1.54
$ws.Range("AG2").Value = 6500
$ws.Range("AH2").Value = 1.61
$ws.Range("AI2").Value = 22.92
$ws.Range("AJ2").Value = 18870000

# Row 3
$ws.Range("D3").Value = 47714
$ws.Range("E3").Value = 6722
$ws.Range("F3").Value = 6722
$ws.Range("G3").Value = 6845
$ws.Range("H3").Value = 5140
$ws.Range("I3").Value = 5094
$ws.Range("J3").Value = 46
$ws.Range("K3").Value = 59079
$ws.Range("L3").Value = 6886
$ws.Range("M3").Value = 52194
$ws.Range("N3").Value = 50773
$ws.Range("O3").Value = 1421
$ws.Range("P3").Value = 944
$ws.Range("Q3").Value = 9449
$ws.Range("R3").Value = -6508
$ws.Range("S3").Value = -2363
$ws.Range("T3").Value = 4799
$ws.Range("U3").Value = 4650
$ws.Range("V3").Value = 529
$ws.Range("W3").Value = 14.09
$ws.Range("X3").Value = 10.77
$ws.Range("Y3").Value = 10.49
$ws.Range("Z3").Value = 8.9
$ws.Range("AA3").Value = 13.19
$ws.Range("AB3").Value = 5216.99
$ws.Range("AC3").Value = 26994
$ws.Range("AD3").Value = 17.37
$ws.Range("AE3").Value = 287269
$ws.Range("AF3").Value = 1.63
$ws.Range("AG3").Value = 8500
$ws.Range("AH3").Value = 1.81
$ws.Range("AI3").Value = 29.49
$ws.Range("AJ3").Value = 18870000

# Row 4
$ws.Range("D4").Value = 58475
$ws.Range("E4").Value = 7647
$ws.Range("F4").Value = 7647
$ws.Range("G4").Value = 8017
$ws.Range("H4").Value = 5946
$ws.Range("I4").Value = 5924
$ws.Range("J4").Value = 22
$ws.Range("K4").Value = 64826
$ws.Range("L4").Value = 8199
$ws.Range("M4").Value = 56627
$ws.Range("N4").Value = 55165
$ws.Range("O4").Value = 1463
$ws.Range("P4").Value = 944
$ws.Range("Q4").Value = 6807
$ws.Range("R4").Value = -4590
$ws.Range("S4").Value = -1565
$ws.Range("T4").Value = 1845
$ws.Range("U4").Value = 4962
$ws.Range("V4").Value = 480
$ws.Range("W4").Value = 13.08
$ws.Range("X4").Value = 10.17
$ws.Range("Y4").Value = 11.18
$ws.Range("Z4").Value = 9.6
$ws.Range("AA4").Value = 14.48
$ws.Range("AB4").Value = 5678.25
$ws.Range("AC4").Value = 31395
$ws.Range("AD4").Value = 15.13
$ws.Range("AE4").Value = 312118
$ws.Range("AF4").Value = 1.52
$ws.Range("AG4").Value = 8500
$ws.Range("AH4").Value = 1.79
$ws.Range("AI4").Value = 25.36
$ws.Range("AJ4").Value = 18870000

# Row 5
$ws.Range("D5").Value = 65967
$ws.Range("E5").Value = 8948
$ws.Range("F5").Value = 8948
$ws.Range("G5").Value = 9257
$ws.Range("H5").Value = 6340
$ws.Range("I5").Value = 6290
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 70381
$ws.Range("L5").Value = 9325
$ws.Range("M5").Value = 61057
$ws.Range("N5").Value = 59603
$ws.Range("O5").Value = 1453
$ws.Range("P5").Value = 944
$ws.Range("Q5").Value = 7968
$ws.Range("R5").Value = -3237
$ws.Range("S5").Value = -1725
$ws.Range("T5").Value = 4970
$ws.Range("U5").Value = 2998
$ws.Range("V5").Value = 338
$ws.Range("W5").Value = 13.56
$ws.Range("X5").Value = 9.609999999999999
$ws.Range("Y5").Value = 10.96
$ws.Range("Z5").Value = 9.380000000000001
$ws.Range("AA5").Value = 15.27
$ws.Range("AB5").Value = 6187.4
$ws.Range("AC5").Value = 33336
$ws.Range("AD5").Value = 14.79
$ws.Range("AE5").Value = 337232
$ws.Range("AF5").Value = 1.46
$ws.Range("AG5").Value = 10000
$ws.Range("AH5").Value = 2.03
$ws.Range("AI5").Value = 28.1
$ws.Range("AJ5").Value = 18870000

# Row 6
$ws.Range("D6").Value = 68833
$ws.Range("E6").Value = 7647
$ws.Range("F6").Value = 7647
$ws.Range("G6").Value = 7687
$ws.Range("H6").Value = 5348
$ws.Range("I6").Value = 5274
$ws.Range("K6").Value = 72259
$ws.Range("L6").Value = 8514
$ws.Range("M6").Value = 63745
$ws.Range("N6").Value = 62500
$ws.Range("P6").Value = 944
$ws.Range("Q6").Value = 8883
$ws.Range("R6").Value = -6497
$ws.Range("S6").Value = -1829
$ws.Range("T6").Value = 3162
$ws.Range("U6").Value = 5722
$ws.Range("V6").Value = 305
$ws.Range("W6").Value = 11.11
$ws.Range("X6").Value = 7.77
$ws.Range("Y6").Value = 8.640000000000001
$ws.Range("Z6").Value = 7.5
$ws.Range("AA6").Value = 13.36
$ws.Range("AB6").Value = 6565.05
$ws.Range("AC6").Value = 27950
$ws.Range("AD6").Value = 15.47
$ws.Range("AE6").Value = 353622
$ws.Range("AF6").Value = 1.22
$ws.Range("AG6").Value = 11000
$ws.Range("AH6").Value = 2.54
$ws.Range("AI6").Value = 36.86
$ws.Range("AJ6").Value = 18870000

# Row 7
$ws.Range("D7").Value = 67691
$ws.Range("E7").Value = 8379
$ws.Range("G7").Value = 9175
$ws.Range("H7").Value = 6547
$ws.Range("I7").Value = 6488
$ws.Range("K7").Value = 77578
$ws.Range("L7").Value = 9248
$ws.Range("M7").Value = 68330
$ws.Range("N7").Value = 67008
$ws.Range("P7").Value = 941
$ws.Range("Q7").Value = 8356
$ws.Range("R7").Value = -4937
$ws.Range("S7").Value = -1914
$ws.Range("T7").Value = 3503
$ws.Range("U7").Value = 5569
$ws.Range("W7").Value = 12.38
$ws.Range("X7").Value = 9.67
$ws.Range("Y7").Value = 10.02
$ws.Range("Z7").Value = 8.74
$ws.Range("AA7").Value = 13.53
$ws.Range("AC7").Value = 34380
$ws.Range("AD7").Value = 11.05
$ws.Range("AE7").Value = 379127
$ws.Range("AF7").Value = 1
$ws.Range("AG7").Value = 11385
$ws.Range("AH7").Value = 3
$ws.Range("AI7").Value = 33.11

# Row 8
$ws.Range("D8").Value = 69362
$ws.Range("E8").Value = 8793
$ws.Range("G8").Value = 9483
$ws.Range("H8").Value = 6784
$ws.Range("I8").Value = 6720
$ws.Range("K8").Value = 82207
$ws.Range("L8").Value = 9209
$ws.Range("M8").Value = 72998
$ws.Range("N8").Value = 71595
$ws.Range("P8").Value = 941
$ws.Range("Q8").Value = 9608
$ws.Range("R8").Value = -5211
$ws.Range("S8").Value = -2034
$ws.Range("T8").Value = 3981
$ws.Range("U8").Value = 5667
$ws.Range("W8").Value = 12.68
$ws.Range("X8").Value = 9.779999999999999
$ws.Range("Y8").Value = 9.699999999999999
$ws.Range("Z8").Value = 8.49
$ws.Range("AA8").Value = 12.62
$ws.Range("AC8").Value = 35610
$ws.Range("AD8").Value = 10.67
$ws.Range("AE8").Value = 405079
$ws.Range("AF8").Value = 0.9399999999999999
$ws.Range("AG8").Value = 11615
$ws.Range("AH8").Value = 3.06
$ws.Range("AI8").Value = 32.62

# Row 9
$ws.Range("D9").Value = 70744
$ws.Range("E9").Value = 9156
$ws.Range("G9").Value = 9894
$ws.Range("H9").Value = 7081
$ws.Range("I9").Value = 7018
$ws.Range("K9").Value = 87294
$ws.Range("L9").Value = 9484
$ws.Range("M9").Value = 77811
$ws.Range("N9").Value = 76324
$ws.Range("P9").Value = 941
$ws.Range("Q9").Value = 9797
$ws.Range("R9").Value = -5405
$ws.Range("S9").Value = -2069
$ws.Range("T9").Value = 3936
$ws.Range("U9").Value = 5811
$ws.Range("W9").Value = 12.94
$ws.Range("X9").Value = 10.01
$ws.Range("Y9").Value = 9.49
$ws.Range("Z9").Value = 8.36
$ws.Range("AA9").Value = 12.19
$ws.Range("AC9").Value = 37191
$ws.Range("AD9").Value = 10.22
$ws.Range("AE9").Value = 431840
$ws.Range("AF9").Value = 0.88
$ws.Range("AG9").Value = 11846
$ws.Range("AH9").Value = 3.12
$ws.Range("AI9").Value = 31.85
